$wb = $excel.ActiveWorkbook

# Update the userName value shared between Register!I2 and LogIn!A2
$wsRegister = $wb.Worksheets.Item("Register")
$wsLogin = $wb.Worksheets.Item("LogIn")

$wsRegister.Range("I2").Value = "newUser34422"
$wsLogin.Range("A2").Value = "newUser34422"

# Select a cell on the Login sheet first to mimic prior state, then move to Register
$wsLogin.Activate()
$wsLogin.Range("A2").Select()

$wsRegister.Activate()
$wsRegister.Range("I6").Select()

$wb.Save()
